$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new weekly blue tank titration row (row 40)
$ws.Range("A40").Value = 20210603
$ws.Range("B40").Value = 2228.9667413407301
$ws.Range("C40").Value = 2224.4699999999998
$ws.Range("D40").Formula = "=100*(B40-C40)/C40"
$ws.Range("E40").Value = 180
$ws.Range("F40").Value = "CRM opened 20210526"

$wb.Save()
